# Bugfixed evaluation and simulated rt_data for components.
# A new (earlier) observation is inserted at the top of the forecast table
# and every y_0_forecast / y_1_forecast value is refreshed with the
# recomputed figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: shift existing data rows 2..18 down to rows 3..19 -------------
# Process bottom-up so we never overwrite a row before reading it.
for ($r = 18; $r -ge 2; $r--) {
  $dest = $r + 1
  for ($col = 1; $col -le 5; $col++) {
    $srcCell = $ws.Cells.Item($r, $col)
    $destCell = $ws.Cells.Item($dest, $col)
    $v = $srcCell.Value2
    if ($v -eq $null) {
      $destCell.ClearContents()
    } else {
      $destCell.Value2 = $v
    }
  }
}

# Row 19 is brand new territory (old sheet only went to row 18) -> make sure
# column A carries the same date style used by every other row in col A.
$ws.Cells.Item(18, 1).Copy()
$ws.Cells.Item(19, 1).PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Step 2: write the new data point into row 2 ----------------------------
$ws.Cells.Item(2, 1).Value2 = 39400
$ws.Cells.Item(2, 2).Value2 = 2007
$ws.Cells.Item(2, 4).Value2 = 2008
$ws.Cells.Item(2, 5).ClearContents()

# --- Step 3: recalculated forecast values (columns C and E) for every row --
$values = @(
  @(2, 5.896808312953783, $null),
  @(3, 7.441962824572235, $null),
  @(4, 6.277541464866987, $null),
  @(5, 6.535114773304773, 6.325696408067327),
  @(6, 5.12051970717502, 4.950888348161886),
  @(7, 3.65682115264816, 3.982564147794321),
  @(8, 2.943878639034381, 4.334309403335435),
  @(9, 1.172679597477866, 2.644356903452572),
  @(10, 2.961845079861303, 3.383932287548697),
  @(11, 2.508469427909898, 3.355044026998955),
  @(12, 3.523703831572056, 3.74984170812418),
  @(13, 1.178605266817589, 2.186196327763934),
  @(14, 3.047037961814492, 2.880436144359444),
  @(15, -0.2228847697281378, 1.982741503124119),
  @(16, -1.165854108406617, 2.782217648649521),
  @(17, 2.501311189006916, 2.985901060752827),
  @(18, 0.6753076481029074, 0.7957830962485257),
  @(19, 2.039329803030121, 2.510359031091491)
)

foreach ($row in $values) {
  $r = $row[0]
  $cVal = $row[1]
  $eVal = $row[2]
  $ws.Cells.Item($r, 3).Value2 = $cVal
  if ($eVal -eq $null) {
    $ws.Cells.Item($r, 5).ClearContents()
  } else {
    $ws.Cells.Item($r, 5).Value2 = $eVal
  }
}
